$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3158036.2
$ws.Range("I4").Value = 129.11765
$ws.Range("J4").Value = 30000248
$ws.Range("K4").Value = 129.11765
$ws.Range("L4").Value = 30000248
$ws.Range("M4").Value = -15.11765
$ws.Range("N4").Value = -30000476
$ws.Range("H12").Value = 537.0769
$ws.Range("I12").Value = 548.375
$ws.Range("K12").Value = 548.375
$ws.Range("M12").Value = -378.375
$ws.Range("H54").Value = 12021
$ws.Range("I54").Value = 7500
$ws.Range("J54").Value = 16542
$ws.Range("K54").Value = 7500
$ws.Range("L54").Value = 16542
$ws.Range("M54").Value = -7014
$ws.Range("N54").Value = -17514
$ws.Range("H132").Value = 2023.5807
$ws.Range("I132").Value = 1860.4445
$ws.Range("K132").Value = 5581.333500000001
$ws.Range("M132").Value = -3051.333500000001
$ws.Range("H137").Value = 1658.7142
$ws.Range("I137").Value = 1379.8889
$ws.Range("J137").Value = 3331.6667
$ws.Range("K137").Value = 4139.6667
$ws.Range("L137").Value = 9995.000100000001
$ws.Range("M137").Value = -1589.6667
$ws.Range("N137").Value = -15095.0001
$ws.Range("H138").Value = 1569713
$ws.Range("I138").Value = 1434.2439
$ws.Range("J138").Value = 3713027.2
$ws.Range("K138").Value = 4302.7317
$ws.Range("L138").Value = 11139081.6
$ws.Range("M138").Value = 837.2682999999997
$ws.Range("N138").Value = -11149361.6
$ws.Range("H141").Value = 4878
$ws.Range("I141").Value = 2940.7307
$ws.Range("K141").Value = 8822.1921
$ws.Range("M141").Value = -3642.1921
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 9999
$ws.Range("J4").Value = 9999
$ws.Range("L4").Value = 9999
$ws.Range("N4").Value = -10231
$ws.Range("H32").Value = 5177.6597
$ws.Range("I32").Value = 5177.6597
$ws.Range("K32").Value = 5177.6597
$ws.Range("M32").Value = -4890.6597
$ws.Range("H45").Value = 1772.9375
$ws.Range("I45").Value = 1364.8334
$ws.Range("K45").Value = 1364.8334
$ws.Range("M45").Value = -987.8334
$ws.Range("H61").Value = 5423.0356
$ws.Range("I61").Value = 4413.8
$ws.Range("J61").Value = 13833.333
$ws.Range("K61").Value = 4413.8
$ws.Range("L61").Value = 13833.333
$ws.Range("M61").Value = -4201.8
$ws.Range("N61").Value = -14257.333
$ws.Range("H74").Value = 3138.1614
$ws.Range("I74").Value = 2989.5
$ws.Range("J74").Value = 3501.5557
$ws.Range("K74").Value = 2989.5
$ws.Range("L74").Value = 3501.5557
$ws.Range("M74").Value = -2115.5
$ws.Range("N74").Value = -5249.5557
$ws.Range("H77").Value = 3138.1614
$ws.Range("I77").Value = 2989.5
$ws.Range("J77").Value = 3501.5557
$ws.Range("K77").Value = 14947.5
$ws.Range("L77").Value = 17507.7785
$ws.Range("M77").Value = -10579.5
$ws.Range("N77").Value = -26243.7785
$ws.Range("H96").Value = 48331.668
$ws.Range("J96").Value = 48331.668
$ws.Range("L96").Value = 48331.668
$ws.Range("N96").Value = -53823.668
$ws.Range("H122").Value = 7410362.5
$ws.Range("I122").Value = 15875075
$ws.Range("J122").Value = 3739
$ws.Range("K122").Value = 47625225
$ws.Range("L122").Value = 11217
$ws.Range("M122").Value = -47622775
$ws.Range("N122").Value = -16117
$ws.Range("H132").Value = 1754.359
$ws.Range("I132").Value = 1687.0541
$ws.Range("K132").Value = 5061.1623
$ws.Range("M132").Value = -2531.1623
$ws.Range("H136").Value = 5423.0356
$ws.Range("I136").Value = 4413.8
$ws.Range("J136").Value = 13833.333
$ws.Range("K136").Value = 13241.4
$ws.Range("L136").Value = 41499.999
$ws.Range("M136").Value = -10691.4
$ws.Range("N136").Value = -46599.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3193.353
$ws.Range("I134").Value = 3733.2222
$ws.Range("J134").Value = 2586
$ws.Range("K134").Value = 11199.6666
$ws.Range("L134").Value = 7758
$ws.Range("M134").Value = -8664.6666
$ws.Range("N134").Value = -12828
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2742.7144
$ws.Range("I2").Value = 1340
$ws.Range("J2").Value = 6249.5
$ws.Range("K2").Value = 1340
$ws.Range("L2").Value = 6249.5
$ws.Range("M2").Value = -1227
$ws.Range("N2").Value = -6475.5
$ws.Range("H22").Value = 2399.2
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1650
$ws.Range("H31").Value = 2944.842
$ws.Range("I31").Value = 1908.75
$ws.Range("J31").Value = 4721
$ws.Range("K31").Value = 1908.75
$ws.Range("L31").Value = 4721
$ws.Range("M31").Value = -1613.75
$ws.Range("N31").Value = -5311
$ws.Range("H34").Value = 2944.842
$ws.Range("I34").Value = 1908.75
$ws.Range("J34").Value = 4721
$ws.Range("K34").Value = 1908.75
$ws.Range("L34").Value = 4721
$ws.Range("M34").Value = -1706.75
$ws.Range("N34").Value = -5125
$ws.Range("H132").Value = 3095.4583
$ws.Range("I132").Value = 1418.0952
$ws.Range("J132").Value = 14837
$ws.Range("K132").Value = 4254.2856
$ws.Range("L132").Value = 44511
$ws.Range("M132").Value = -1724.2856
$ws.Range("N132").Value = -49571
$ws.Range("H134").Value = 24649.072
$ws.Range("I134").Value = 27216.96
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 81650.88
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -79115.88
$ws.Range("N134").Value = -14820
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1747.0769
$ws.Range("I26").Value = 1701.3636
$ws.Range("J26").Value = 1998.5
$ws.Range("K26").Value = 5104.0908
$ws.Range("L26").Value = 5995.5
$ws.Range("M26").Value = -4816.0908
$ws.Range("N26").Value = -6571.5
$ws.Range("H34").Value = 3739.0715
$ws.Range("J34").Value = 3852.6758
$ws.Range("L34").Value = 11558.0274
$ws.Range("N34").Value = -11726.0274
$ws.Range("H55").Value = 441966.38
$ws.Range("I55").Value = 720174.7
$ws.Range("J55").Value = 117390
$ws.Range("K55").Value = 2160524.1
$ws.Range("L55").Value = 352170
$ws.Range("M55").Value = -2160347.1
$ws.Range("N55").Value = -352524
$ws.Range("H98").Value = 1097.7142
$ws.Range("I98").Value = 146.66667
$ws.Range("J98").Value = 1811
$ws.Range("K98").Value = 440.00001
$ws.Range("L98").Value = 5433
$ws.Range("M98").Value = 1057.99999
$ws.Range("N98").Value = -8429
$ws.Range("H122").Value = 2343.4
$ws.Range("J122").Value = 2439.4285
$ws.Range("L122").Value = 21954.8565
$ws.Range("N122").Value = -26854.8565
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 4500
$ws.Range("K5").Value = 4500
$ws.Range("M5").Value = -4388
$ws.Range("H58").Value = 13000.333
$ws.Range("I58").Value = 14500
$ws.Range("K58").Value = 14500
$ws.Range("M58").Value = -14223
$ws.Range("H92").Value = 18062.25
$ws.Range("J92").Value = 18062.25
$ws.Range("L92").Value = 18062.25
$ws.Range("N92").Value = -21806.25
$ws.Range("H113").Value = 4410.6924
$ws.Range("J113").Value = 6299.2173
$ws.Range("L113").Value = 6299.2173
$ws.Range("N113").Value = -10639.2173
$ws.Range("H122").Value = 4567
$ws.Range("I122").Value = 2591.1428
$ws.Range("K122").Value = 7773.428400000001
$ws.Range("M122").Value = -5323.428400000001
$ws.Range("H123").Value = 33333
$ws.Range("J123").Value = 33333
$ws.Range("L123").Value = 33333
$ws.Range("N123").Value = -38233
$ws.Range("H132").Value = 2954.2
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2954.2
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8862.599999999999
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13922.6
$ws.Range("H138").Value = 110143
$ws.Range("J138").Value = 110143
$ws.Range("L138").Value = 110143
$ws.Range("N138").Value = -120423
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 919.5714
$ws.Range("J46").Value = 1087.6
$ws.Range("L46").Value = 1087.6
$ws.Range("N46").Value = -1463.6
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H122").Value = 6428.0586
$ws.Range("I122").Value = 5031.5557
$ws.Range("K122").Value = 15094.6671
$ws.Range("M122").Value = -12644.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25799
$ws.Range("I54").Value = 14999
$ws.Range("K54").Value = 14999
$ws.Range("M54").Value = -14479
$ws.Range("H132").Value = 2054.449
$ws.Range("I132").Value = 1803.9736
$ws.Range("J132").Value = 2919.7273
$ws.Range("K132").Value = 5411.9208
$ws.Range("L132").Value = 8759.1819
$ws.Range("M132").Value = -2881.9208
$ws.Range("N132").Value = -13819.1819
$ws.Range("H136").Value = 2944.6191
$ws.Range("I136").Value = 792.93335
$ws.Range("K136").Value = 2378.80005
$ws.Range("M136").Value = 171.1999500000002
